# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Gungnir_Profits workbook
# (values refreshed by the scheduled pricing runner)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2562.162
$ws.Range("I64").Value = 2540
$ws.Range("K64").Value = 2540
$ws.Range("M64").Value = -2292
$ws.Range("H67").Value = 2562.162
$ws.Range("I67").Value = 2540
$ws.Range("K67").Value = 2540
$ws.Range("M67").Value = -1682
$ws.Range("H70").Value = 3034.1667
$ws.Range("I70").Value = 3551.25
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 10653.75
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -10383.75
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 3034.1667
$ws.Range("I73").Value = 3551.25
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 10653.75
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -9717.75
$ws.Range("N73").Value = -7872
$ws.Range("H82").Value = 996.6667
$ws.Range("I82").Value = 996.6667
$ws.Range("K82").Value = 2990.0001
$ws.Range("M82").Value = -2584.0001
$ws.Range("H85").Value = 996.6667
$ws.Range("I85").Value = 996.6667
$ws.Range("K85").Value = 2990.0001
$ws.Range("M85").Value = -1586.0001
$ws.Range("H100").Value = 14000
$ws.Range("I100").Value = 16857.143
$ws.Range("J100").Value = 9000
$ws.Range("K100").Value = 16857.143
$ws.Range("L100").Value = 9000
$ws.Range("M100").Value = -16316.143
$ws.Range("N100").Value = -10082
$ws.Range("H141").Value = 5517
$ws.Range("I141").Value = 2286.889
$ws.Range("J141").Value = 20052.5
$ws.Range("K141").Value = 6860.667
$ws.Range("L141").Value = 60157.5
$ws.Range("M141").Value = -1680.667
$ws.Range("N141").Value = -70517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19610912
$ws.Range("I32").Value = 2880.7273
$ws.Range("K32").Value = 2880.7273
$ws.Range("M32").Value = -2593.7273
$ws.Range("H45").Value = 38247.555
$ws.Range("I45").Value = 84059
$ws.Range("J45").Value = 1598.4
$ws.Range("K45").Value = 84059
$ws.Range("L45").Value = 1598.4
$ws.Range("M45").Value = -83682
$ws.Range("N45").Value = -2352.4
$ws.Range("H97").Value = 2099.9
$ws.Range("I97").Value = 2236.5
$ws.Range("J97").Value = 1895
$ws.Range("K97").Value = 2236.5
$ws.Range("L97").Value = 1895
$ws.Range("M97").Value = -1740.5
$ws.Range("N97").Value = -2887
$ws.Range("H132").Value = 1279913.5
$ws.Range("I132").Value = 984.1053000000001
$ws.Range("J132").Value = 7354828
$ws.Range("K132").Value = 2952.3159
$ws.Range("L132").Value = 22064484
$ws.Range("M132").Value = -422.3159000000001
$ws.Range("N132").Value = -22069544

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 930.2
$ws.Range("I94").Value = 912.75
$ws.Range("K94").Value = 912.75
$ws.Range("M94").Value = -461.75
$ws.Range("H99").Value = 1935.3158
$ws.Range("I99").Value = 1197.25
$ws.Range("K99").Value = 1197.25
$ws.Range("M99").Value = 300.75
$ws.Range("H105").Value = 47620850
$ws.Range("I105").Value = 1721.2667
$ws.Range("J105").Value = 166668670
$ws.Range("K105").Value = 1721.2667
$ws.Range("L105").Value = 166668670
$ws.Range("M105").Value = 25.7333000000001
$ws.Range("N105").Value = -166672164
$ws.Range("H134").Value = 2139958.2
$ws.Range("I134").Value = 1175.5405
$ws.Range("J134").Value = 7415622
$ws.Range("K134").Value = 3526.6215
$ws.Range("L134").Value = 22246866
$ws.Range("M134").Value = -991.6215000000002
$ws.Range("N134").Value = -22251936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1856386
$ws.Range("I31").Value = 2227094.8
$ws.Range("K31").Value = 2227094.8
$ws.Range("M31").Value = -2226799.8
$ws.Range("H34").Value = 1856386
$ws.Range("I34").Value = 2227094.8
$ws.Range("K34").Value = 2227094.8
$ws.Range("M34").Value = -2226892.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 9631.272000000001
$ws.Range("I68").Value = 461.75
$ws.Range("J68").Value = 34083.332
$ws.Range("K68").Value = 1385.25
$ws.Range("L68").Value = 102249.996
$ws.Range("M68").Value = -574.25
$ws.Range("N68").Value = -103871.996
$ws.Range("H71").Value = 9631.272000000001
$ws.Range("I71").Value = 461.75
$ws.Range("J71").Value = 34083.332
$ws.Range("K71").Value = 4155.75
$ws.Range("L71").Value = 306749.988
$ws.Range("M71").Value = -99.75
$ws.Range("N71").Value = -314861.988
$ws.Range("H94").Value = 2620
$ws.Range("I94").Value = 1500
$ws.Range("J94").Value = 2900
$ws.Range("K94").Value = 4500
$ws.Range("L94").Value = 8700
$ws.Range("M94").Value = -3824
$ws.Range("N94").Value = -10052
$ws.Range("H97").Value = 3665.5
$ws.Range("I97").Value = 5000
$ws.Range("J97").Value = 3398.6
$ws.Range("K97").Value = 15000
$ws.Range("L97").Value = 10195.8
$ws.Range("M97").Value = -14504
$ws.Range("N97").Value = -11187.8
$ws.Range("H100").Value = 2000
$ws.Range("J100").Value = 2000
$ws.Range("L100").Value = 6000
$ws.Range("N100").Value = -7622
$ws.Range("H129").Value = 1392.4117
$ws.Range("I129").Value = 1230
$ws.Range("J129").Value = 1402.5625
$ws.Range("K129").Value = 3690
$ws.Range("L129").Value = 4207.6875
$ws.Range("M129").Value = 1310
$ws.Range("N129").Value = -14207.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 41666.668
$ws.Range("J27").Value = 41666.668
$ws.Range("L27").Value = 41666.668
$ws.Range("N27").Value = -41998.668
$ws.Range("H97").Value = 1548.0834
$ws.Range("I97").Value = 1508.8182
$ws.Range("J97").Value = 1980
$ws.Range("K97").Value = 1508.8182
$ws.Range("L97").Value = 1980
$ws.Range("M97").Value = -1012.8182
$ws.Range("N97").Value = -2972
$ws.Range("H111").Value = 18096.666
$ws.Range("J111").Value = 18096.666
$ws.Range("L111").Value = 18096.666
$ws.Range("N111").Value = -24230.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5600811
$ws.Range("I22").Value = 14932113
$ws.Range("J22").Value = 2029.7
$ws.Range("K22").Value = 14932113
$ws.Range("L22").Value = 2029.7
$ws.Range("M22").Value = -14931818
$ws.Range("N22").Value = -2619.7
$ws.Range("H27").Value = 5600811
$ws.Range("I27").Value = 14932113
$ws.Range("J27").Value = 2029.7
$ws.Range("K27").Value = 14932113
$ws.Range("L27").Value = 2029.7
$ws.Range("M27").Value = -14932006
$ws.Range("N27").Value = -2243.7
$ws.Range("H68").Value = 6945771
$ws.Range("I68").Value = 6945771
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 6945771
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -6945022
$ws.Range("N68").ClearContents()
$ws.Range("H70").Value = 30160
$ws.Range("J70").Value = 30160
$ws.Range("L70").Value = 30160
$ws.Range("N70").Value = -30700
$ws.Range("H71").Value = 6945771
$ws.Range("I71").Value = 6945771
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 34728855
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -34725111
$ws.Range("N71").ClearContents()
$ws.Range("H73").Value = 30160
$ws.Range("J73").Value = 30160
$ws.Range("L73").Value = 30160
$ws.Range("N73").Value = -32032
$ws.Range("H75").Value = 47780
$ws.Range("J75").Value = 47780
$ws.Range("L75").Value = 47780
$ws.Range("N75").Value = -49652
$ws.Range("H78").Value = 47780
$ws.Range("J78").Value = 47780
$ws.Range("L78").Value = 143340
$ws.Range("N78").Value = -152700
$ws.Range("H82").Value = 71429320
$ws.Range("I82").Value = 71429320
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 71429320
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -71428959
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 71429320
$ws.Range("I85").Value = 71429320
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 71429320
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -71428072
$ws.Range("N85").ClearContents()
$ws.Range("H93").Value = 1131.2
$ws.Range("I93").Value = 1187.4286
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 1187.4286
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 60.57140000000004
$ws.Range("N93").Value = -3496
$ws.Range("H100").Value = 2605.8125
$ws.Range("I100").Value = 3200.6
$ws.Range("J100").Value = 2495.6667
$ws.Range("K100").Value = 3200.6
$ws.Range("L100").Value = 2495.6667
$ws.Range("M100").Value = -2659.6
$ws.Range("N100").Value = -3577.6667
$ws.Range("H110").Value = 26911
$ws.Range("J110").Value = 26911
$ws.Range("L110").Value = 26911
$ws.Range("N110").Value = -35091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1068.909
$ws.Range("I81").Value = 861.55554
$ws.Range("J81").Value = 2002
$ws.Range("K81").Value = 1723.11108
$ws.Range("L81").Value = 4004
$ws.Range("M81").Value = -662.1110799999999
$ws.Range("N81").Value = -6126
$ws.Range("H84").Value = 1068.909
$ws.Range("I84").Value = 861.55554
$ws.Range("J84").Value = 2002
$ws.Range("K84").Value = 8615.555399999999
$ws.Range("L84").Value = 20020
$ws.Range("M84").Value = -3311.555399999999
$ws.Range("N84").Value = -30628
$ws.Range("H96").Value = 1922.2222
$ws.Range("J96").Value = 3500
$ws.Range("L96").Value = 3500
$ws.Range("N96").Value = -6246
